$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 616.3
$ws.Range("K5").Value = 437
$ws.Range("K6").Value = 746.8
